# Insert a new weekly price record at row 124 in the "Hortaliza, Terminal La
# Palmera de La Serena - Zanahoria" sheet. This pushes the existing rows
# 124:188 down to 125:189 (preserving all their data), and the newly
# created row 124 is populated with a new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 124, shifting rows 124:188 -> 125:189.
$ws.Rows.Item(124).Insert()

# Populate the new row 124 with the new weekly record.
$ws.Range("A124").Value = 8
$ws.Range("B124").Value = "Terminal La Palmera de La Serena"
$ws.Range("C124").Value = "Coquimbo"
$ws.Range("D124").Value = 44455
$ws.Range("E124").Value = 4
$ws.Range("F124").Value = 100114013
$ws.Range("G124").Value = "Zanahoria"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 720
$ws.Range("K124").Value = 4500
$ws.Range("L124").Value = 5000
$ws.Range("M124").Value = 4750
$ws.Range("N124").Value = "$/saco 20 kilos"
$ws.Range("O124").Value = "Provincia del Elquí"
$ws.Range("P124").Value = 238
$ws.Range("Q124").Value = 20
$ws.Range("R124").Value = "Hortaliza"
